$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency price / 1h-volume-change data.
# Rows 24/25 and 37/38 also swap coin identity (re-ranked by the scraper).
#
# Price (column D) values are stored as literal text in the source sheet
# (e.g. "43.446.43", "0.810"), not numbers -- a leading apostrophe forces
# Excel to keep them as text instead of auto-coercing them to numbers
# (which would silently drop meaningful trailing zeros, e.g. "0.810" -> 0.81).

# Row 2: update D, E
$ws.Range("D2").Value = "43.446.43"
$ws.Range("E2").Value = "  +2.85%  "

# Row 3: update D, E
$ws.Range("D3").Value = "2.311.55"
$ws.Range("E3").Value = "  +1.90%  "

# Row 4: update E
$ws.Range("E4").Value = "  +0.01%  "

# Row 5: update D, E
$ws.Range("D5").Value = "'311.36"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6: update D, E
$ws.Range("D6").Value = "'102.09"
$ws.Range("E6").Value = "  +4.99%  "

# Row 7: update E
$ws.Range("E7").Value = "  +1.39%  "

# Row 8: update E
$ws.Range("E8").Value = "  +0.01%  "

# Row 9: update E
$ws.Range("E9").Value = "  +7.45%  "

# Row 10: update D, E
$ws.Range("D10").Value = "'35.82"
$ws.Range("E10").Value = "  +1.85%  "

# Row 11: update E
$ws.Range("E11").Value = "  +2.94%  "

# Row 12: update E
$ws.Range("E12").Value = "  -0.71%  "

# Row 13: update E
$ws.Range("E13").Value = "  +0.75%  "

# Row 14: update D, E
$ws.Range("D14").Value = "2.669.02"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15: update D, E
$ws.Range("D15").Value = "'14.97"
$ws.Range("E15").Value = "  +1.18%  "

# Row 16: update D, E
$ws.Range("D16").Value = "2.308.41"
$ws.Range("E16").Value = "  +2.36%  "

# Row 17: update D, E
$ws.Range("D17").Value = "'0.810"
$ws.Range("E17").Value = "  +2.22%  "

# Row 18: update D, E
$ws.Range("D18").Value = "43.348.68"
$ws.Range("E18").Value = "  +2.93%  "

# Row 19: update D, E
$ws.Range("D19").Value = "'12.35"
$ws.Range("E19").Value = "  +0.44%  "

# Row 20: update D, E
$ws.Range("D20").Value = "0.0₃0931"
$ws.Range("E20").Value = "  +2.75%  "

# Row 21: update D, E
$ws.Range("D21").Value = "'6.17"
$ws.Range("E21").Value = "  +2.26%  "

# Row 22: update D, E
$ws.Range("D22").Value = "'68.09"
$ws.Range("E22").Value = "  +0.24%  "

# Row 23: update D, E
$ws.Range("D23").Value = "'241.62"
$ws.Range("E23").Value = "  +1.52%  "

# Row 24: update B, C, D, E
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = "  +2.14%  "

# Row 25: update B, C, D, E
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'2.01"
$ws.Range("E25").Value = "  +2.29%  "

# Row 26: update E
$ws.Range("E26").Value = "  +0.12%  "

# Row 27: update D, E
$ws.Range("D27").Value = "'3.99"
$ws.Range("E27").Value = "  -1.63%  "

# Row 28: update D, E
$ws.Range("D28").Value = "'24.62"
$ws.Range("E28").Value = "  +4.49%  "

# Row 29: update D, E
$ws.Range("D29").Value = "'36.76"
$ws.Range("E29").Value = "  -2.93%  "

# Row 30: update D, E
$ws.Range("D30").Value = "'9.65"
$ws.Range("E30").Value = "  +1.04%  "

# Row 31: update E
$ws.Range("E31").Value = "  +0.27%  "

# Row 32: update D, E
$ws.Range("D32").Value = "'167.47"
$ws.Range("E32").Value = "  +3.26%  "

# Row 33: update D, E
$ws.Range("D33").Value = "'5.29"
$ws.Range("E33").Value = "  +0.79%  "

# Row 34: update E
$ws.Range("E34").Value = "  +0.04%  "

# Row 35: update E
$ws.Range("E35").Value = "  +0.83%  "

# Row 36: update E
$ws.Range("E36").Value = "  +5.66%  "

# Row 37: update B, C, D, E
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'17.66"
$ws.Range("E37").Value = "  -0.26%  "

# Row 38: update B, C, D, E
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  -2.83%  "

# Row 39: update D, E
$ws.Range("D39").Value = "'1.90"
$ws.Range("E39").Value = "  +4.02%  "

# Row 40: update E
$ws.Range("E40").Value = "  +1.32%  "

# Row 41: update E
$ws.Range("E41").Value = "  +1.51%  "

# Row 42: update D, E
$ws.Range("D42").Value = "'4.36"
$ws.Range("E42").Value = "  +7.14%  "

# Row 43: update D, E
$ws.Range("D43").Value = "'2.31"
$ws.Range("E43").Value = "  -0.87%  "

# Row 44: update E
$ws.Range("E44").Value = "  +2.86%  "

# Row 45: update D, E
$ws.Range("D45").Value = "1.972.92"
$ws.Range("E45").Value = "  +1.21%  "

# Row 46: update E
$ws.Range("E46").Value = "  -1.00%  "

# Row 47: update D, E
$ws.Range("D47").Value = "'2.98"
$ws.Range("E47").Value = "  +2.06%  "

# Row 48: update E
$ws.Range("E48").Value = "  +0.49%  "

# Row 49: update D, E
$ws.Range("D49").Value = "'55.73"
$ws.Range("E49").Value = "  +3.93%  "

# Row 50: update E
$ws.Range("E50").Value = "  +3.26%  "

# Row 51: update E
$ws.Range("E51").Value = "  +7.15%  "
